$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A (closest achievable value to target XML width 147.54296875)
$ws.Columns.Item(1).ColumnWidth = 146.7

# New rows of data (rows 31-35)
$newData = @(
    @("tf-idf, stop words, title 0.6, text 0.3, anchor 0.1  no disambiguation pages, collect all docs of query in one list", 0.23043, 21.61),
    @("tf-idf, stop words, title 0.6, text 0.3, anchor 0.1  +pagerank 0.5 no disambiguation pages, collect all docs of query in one list, query expansion 2 no duplicates", 0.1806, 38.599),
    @("tf-idf, stop words, title 0.7, text 0.2, anchor 0.1   +pagerank 0.2 no disambiguation pages, with anchor disambiguation, collect all docs of query in one list", 0.21176, 18.739),
    @("tf-idf, stop words, title 0.7, text 0.2, anchor 0.1   +pagerank 0.2 no disambiguation pages, with anchor disambiguation, collect all docs of query in one list, skip short docs (>0.8)", 0.2254, 9.7047),
    @("tf-idf, stop words, title 0.6, text 0.3, anchor 0.1   +pagerank 0.4 no disambiguation pages, with anchor disambiguation, collect all docs of query in one list, skip short docs (>0.8)", 0.2437, 20.9025)
)

$row = 31
foreach ($entry in $newData) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

$ws.Range("A29").Select()
